{"js": "// Replace the math-drill answers in the single table of the document.\n// Each table cell holds exactly one equation (e.g. \"8+5=13\"). The commit\n// swaps every equation for a new one, keyed strictly by its position\n// (row, column) in the table -- NOT by its old text -- because a couple\n// of old equations repeat verbatim (\"38+7=45\" appears twice) while\n// mapping to different replacements, so a text-based find/replace would\n// be ambiguous. Position-based replacement is therefore the only\n// reliable technique.\n\n// New answers, row-major, 20 rows x 5 columns (100 cells total).\nconst NEW_VALUES = [\n  [\"57-8=49\", \"41-4=37\", \"24+69=93\", \"27+29=56\", \"85-9=76\"],\n  [\"19+22=41\", \"40-21=19\", \"65-58=7\", \"80-52=28\", \"84-46=38\"],\n  [\"18+29=47\", \"73-14=59\", \"72-14=58\", \"66+18=84\", \"61-49=12\"],\n  [\"68+27=95\", \"7+78=85\", \"76-57=19\", \"57+18=75\", \"14+68=82\"],\n  [\"41-25=16\", \"96-77=19\", \"28+33=61\", \"44+8=52\", \"27+6=33\"],\n  [\"64+9=73\", \"6+55=61\", \"91-88=3\", \"33+28=61\", \"8+73=81\"],\n  [\"43+28=71\", \"9+55=64\", \"97-68=29\", \"48+5=53\", \"12+39=51\"],\n  [\"64-39=25\", \"53+38=91\", \"67-9=58\", \"33+59=92\", \"92-27=65\"],\n  [\"12+49=61\", \"95-88=7\", \"40-38=2\", \"81-2=79\", \"18+76=94\"],\n  [\"53+29=82\", \"47-8=39\", \"87-9=78\", \"44-37=7\", \"48+29=77\"],\n  [\"55-7=48\", \"8+79=87\", \"72-8=64\", \"90-35=55\", \"73-9=64\"],\n  [\"98-19=79\", \"70-8=62\", \"62-13=49\", \"50-11=39\", \"80-3=77\"],\n  [\"84-38=46\", \"22-17=5\", \"38+46=84\", \"92-24=68\", \"72-23=49\"],\n  [\"17+64=81\", \"70-17=53\", \"4+39=43\", \"51-22=29\", \"93-38=55\"],\n  [\"87-19=68\", \"81-4=77\", \"75-29=46\", \"4+58=62\", \"35+47=82\"],\n  [\"9+86=95\", \"64+18=82\", \"38+45=83\", \"9+57=66\", \"2+49=51\"],\n  [\"16+19=35\", \"48+44=92\", \"43+38=81\", \"56+35=91\", \"46+7=53\"],\n  [\"64+18=82\", \"82-35=47\", \"19+76=95\", \"38+13=51\", \"97-59=38\"],\n  [\"8+18=26\", \"25+19=44\", \"88+7=95\", \"24+7=31\", \"89+2=91\"],\n  [\"48+8=56\", \"19+17=36\", \"94-35=59\", \"45+29=74\", \"14+67=81\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Apply the whole grid in one shot (mirrors Table.values setter semantics).\ntable.values = NEW_VALUES;\nawait context.sync();\n", "ps1": "# Replace the math-drill answers in the single table of the document.\n# Each table cell holds exactly one equation (e.g. \"8+5=13\"). The commit\n# swaps every equation for a new one, keyed strictly by its position\n# (row, column) in the table -- NOT by its old text -- because a couple\n# of old equations repeat verbatim (\"38+7=45\" appears twice) while\n# mapping to different replacements, so a text-based find/replace would\n# be ambiguous. Position-based replacement is therefore the only\n# reliable technique.\n\n# New answers, row-major, 20 rows x 5 columns (100 cells total).\n$NEW_VALUES = @(\n    @(\"57-8=49\", \"41-4=37\", \"24+69=93\", \"27+29=56\", \"85-9=76\"),\n    @(\"19+22=41\", \"40-21=19\", \"65-58=7\", \"80-52=28\", \"84-46=38\"),\n    @(\"18+29=47\", \"73-14=59\", \"72-14=58\", \"66+18=84\", \"61-49=12\"),\n    @(\"68+27=95\", \"7+78=85\", \"76-57=19\", \"57+18=75\", \"14+68=82\"),\n    @(\"41-25=16\", \"96-77=19\", \"28+33=61\", \"44+8=52\", \"27+6=33\"),\n    @(\"64+9=73\", \"6+55=61\", \"91-88=3\", \"33+28=61\", \"8+73=81\"),\n    @(\"43+28=71\", \"9+55=64\", \"97-68=29\", \"48+5=53\", \"12+39=51\"),\n    @(\"64-39=25\", \"53+38=91\", \"67-9=58\", \"33+59=92\", \"92-27=65\"),\n    @(\"12+49=61\", \"95-88=7\", \"40-38=2\", \"81-2=79\", \"18+76=94\"),\n    @(\"53+29=82\", \"47-8=39\", \"87-9=78\", \"44-37=7\", \"48+29=77\"),\n    @(\"55-7=48\", \"8+79=87\", \"72-8=64\", \"90-35=55\", \"73-9=64\"),\n    @(\"98-19=79\", \"70-8=62\", \"62-13=49\", \"50-11=39\", \"80-3=77\"),\n    @(\"84-38=46\", \"22-17=5\", \"38+46=84\", \"92-24=68\", \"72-23=49\"),\n    @(\"17+64=81\", \"70-17=53\", \"4+39=43\", \"51-22=29\", \"93-38=55\"),\n    @(\"87-19=68\", \"81-4=77\", \"75-29=46\", \"4+58=62\", \"35+47=82\"),\n    @(\"9+86=95\", \"64+18=82\", \"38+45=83\", \"9+57=66\", \"2+49=51\"),\n    @(\"16+19=35\", \"48+44=92\", \"43+38=81\", \"56+35=91\", \"46+7=53\"),\n    @(\"64+18=82\", \"82-35=47\", \"19+76=95\", \"38+13=51\", \"97-59=38\"),\n    @(\"8+18=26\", \"25+19=44\", \"88+7=95\", \"24+7=31\", \"89+2=91\"),\n    @(\"48+8=56\", \"19+17=36\", \"94-35=59\", \"45+29=74\", \"14+67=81\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $NEW_VALUES.Count; $r++) {\n    $rowVals = $NEW_VALUES[$r - 1]\n    for ($c = 1; $c -le $rowVals.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowVals[$c - 1]\n    }\n}\n"}
